$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 320-321, pushing the existing data
# (currently rows 320-341) down to rows 322-343.
$ws.Range("A320:A321").EntireRow.Insert()

# New row 320: Alcachofa, Argentina(o), Primera
$ws.Cells.Item(320, 1).Value = 3
$ws.Cells.Item(320, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(320, 3).Value = "Coquimbo"
$ws.Cells.Item(320, 4).Value = 44753
$ws.Cells.Item(320, 5).Value = 5
$ws.Cells.Item(320, 6).Value = 100112013
$ws.Cells.Item(320, 7).Value = "Alcachofa"
$ws.Cells.Item(320, 8).Value = "Argentina(o)"
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 115
$ws.Cells.Item(320, 11).Value = 14000
$ws.Cells.Item(320, 12).Value = 15000
$ws.Cells.Item(320, 13).Value = 14478
$ws.Cells.Item(320, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(320, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(320, 16).Value = 290
$ws.Cells.Item(320, 17).Value = 50
$ws.Cells.Item(320, 18).Value = "Hortaliza"

# New row 321: Alcachofa, Española, Extra
$ws.Cells.Item(321, 1).Value = 3
$ws.Cells.Item(321, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(321, 3).Value = "Coquimbo"
$ws.Cells.Item(321, 4).Value = 44753
$ws.Cells.Item(321, 5).Value = 5
$ws.Cells.Item(321, 6).Value = 100112013
$ws.Cells.Item(321, 7).Value = "Alcachofa"
$ws.Cells.Item(321, 8).Value = "Española"
$ws.Cells.Item(321, 9).Value = "Extra"
$ws.Cells.Item(321, 10).Value = 125
$ws.Cells.Item(321, 11).Value = 17000
$ws.Cells.Item(321, 12).Value = 18000
$ws.Cells.Item(321, 13).Value = 17480
$ws.Cells.Item(321, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(321, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(321, 16).Value = 583
$ws.Cells.Item(321, 17).Value = 30
$ws.Cells.Item(321, 18).Value = "Hortaliza"
